$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new "Tasks vs Requirements" traceability table ---
# Write order matters: it determines the shared-string table index
# assignment order, so it is chosen to reproduce the target file.
$ws.Range("F20").Value = "Tasks"
$ws.Range("F19").Value = "Requirements"

$ws.Range("H20").Value = "R1"
$ws.Range("I20").Value = "R2"
$ws.Range("J20").Value = "R3"
$ws.Range("K20").Value = "R4"
$ws.Range("L20").Value = "R5"

$ws.Range("G21").Value = "T1"
$ws.Range("G22").Value = "T2"
$ws.Range("G23").Value = "T3"
$ws.Range("G24").Value = "T4"
$ws.Range("G25").Value = "T5"

$ws.Range("M20").Value = "R6"

$ws.Range("H21").Value = "x"
$ws.Range("I21").Value = "x"
$ws.Range("K21").Value = "x"
$ws.Range("M22").Value = "x"
$ws.Range("L23").Value = "x"
$ws.Range("J25").Value = "x"

# --- Merge the header cells ---
$ws.Range("F19:M19").Merge()
$ws.Range("F20:F25").Merge()

# --- Alignment (set before borders so no unused "bordered, unaligned"
#     intermediate styles are left behind in the style table) ---
$r1 = $ws.Range("F19:M19")
$r1.HorizontalAlignment = -4108

$r2 = $ws.Range("F20:F25")
$r2.HorizontalAlignment = -4108
$r2.VerticalAlignment = -4108

# --- Borders ---
# "Requirements" merged banner row: outline box around the merged range,
# applied edge-by-edge (top/bottom first, then left/right) so each cell
# in the range keeps only the border segments touching the outside edge.
$r1.Borders.Item(8).LineStyle = 1
$r1.Borders.Item(9).LineStyle = 1
$r1.Borders.Item(7).LineStyle = 1
$r1.Borders.Item(10).LineStyle = 1

# "Tasks" merged column: same outline treatment.
$r2.Borders.Item(8).LineStyle = 1
$r2.Borders.Item(9).LineStyle = 1
$r2.Borders.Item(7).LineStyle = 1
$r2.Borders.Item(10).LineStyle = 1

# Matrix body: every cell boxed individually (All Borders).
$ws.Range("G20:M25").Borders.LineStyle = 1

# --- Selection matches the saved file ---
$ws.Range("O25").Select()
